# CI: Update Excel counters (state_counters + packages)
# Bump the SPA_Last counter for state "MD" (row 2) from 9371 to 9372.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StateCounters")

$ws.Range("B2").Value = 9372
